# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update the "K" column (column G) with newly calculated strikeout values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 0
    3  = 0
    4  = 0
    5  = 1
    6  = 0
    7  = 0
    8  = 1
    9  = 3
    10 = 4
    11 = 0
    12 = 1
    13 = 0
    14 = 1
    15 = 0
    16 = 1
    17 = 1
    18 = 1
    19 = 1
    20 = 1
    21 = 2
    22 = 1
    23 = 0
    24 = 3
    25 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
